$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 999.3333
$ws.Range("I101").Value = 999.3333
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 2997.9999
$ws.Range("L101").Value = 0
$ws.Range("M101").Value = -1375.9999
$ws.Range("N101").ClearContents()

$ws.Range("H111").Value = 750
$ws.Range("I111").Value = 750
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 2250
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = 817
$ws.Range("N111").ClearContents()

$ws.Range("H113").Value = 2319.25
$ws.Range("I113").Value = 2388.75
$ws.Range("J113").Value = 2249.75
$ws.Range("K113").Value = 2388.75
$ws.Range("L113").Value = 2249.75
$ws.Range("M113").Value = 865.25
$ws.Range("N113").Value = -8757.75

$ws.Range("H116").Value = 3000
$ws.Range("I116").Value = 3000
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 3000
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 442

$ws.Range("H132").Value = 9499.333000000001
$ws.Range("I132").Value = 11999
$ws.Range("J132").Value = 4500
$ws.Range("K132").Value = 35997
$ws.Range("L132").Value = 13500
$ws.Range("M132").Value = -33467
$ws.Range("N132").Value = -18560

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H35").Value = 1500
$ws.Range("I35").Value = 1500
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 1500
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -1094

$ws.Range("H74").Value = 4088.6667
$ws.Range("I74").Value = 1484.6364
$ws.Range("J74").Value = 11249.75
$ws.Range("K74").Value = 1484.6364
$ws.Range("L74").Value = 11249.75
$ws.Range("M74").Value = -610.6364000000001

$ws.Range("H77").Value = 4088.6667
$ws.Range("I77").Value = 1484.6364
$ws.Range("J77").Value = 11249.75
$ws.Range("K77").Value = 7423.182000000001
$ws.Range("L77").Value = 56248.75
$ws.Range("M77").Value = -3055.182000000001

$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("M102").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 190.82222
$ws.Range("I7").Value = 208.33333
$ws.Range("J7").Value = 184.45454
$ws.Range("K7").Value = 208.33333
$ws.Range("L7").Value = 184.45454
$ws.Range("M7").Value = -95.33332999999999
$ws.Range("N7").Value = -410.45454

$ws.Range("H22").Value = 195
$ws.Range("I22").Value = 217
$ws.Range("J22").Value = 129
$ws.Range("K22").Value = 217
$ws.Range("L22").Value = 129
$ws.Range("M22").Value = 133
$ws.Range("N22").Value = -829

$ws.Range("H31").Value = 9796.6
$ws.Range("I31").Value = 5665.5557
$ws.Range("J31").Value = 15993.167
$ws.Range("K31").Value = 5665.5557
$ws.Range("L31").Value = 15993.167
$ws.Range("M31").Value = -5370.5557

$ws.Range("H34").Value = 9796.6
$ws.Range("I34").Value = 5665.5557
$ws.Range("J34").Value = 15993.167
$ws.Range("K34").Value = 5665.5557
$ws.Range("L34").Value = 15993.167
$ws.Range("M34").Value = -5463.5557

$ws.Range("H58").Value = 13608
$ws.Range("I58").Value = 9549.333000000001
$ws.Range("J58").Value = 17666.666
$ws.Range("K58").Value = 9549.333000000001
$ws.Range("L58").Value = 17666.666
$ws.Range("M58").Value = -9346.333000000001

$ws.Range("H132").Value = 7000
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 7000
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 21000
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -26060

$ws.Range("H136").Value = 13608
$ws.Range("I136").Value = 9549.333000000001
$ws.Range("J136").Value = 17666.666
$ws.Range("K136").Value = 28647.999
$ws.Range("L136").Value = 52999.99800000001
$ws.Range("M136").Value = -26097.999

$ws.Range("H141").Value = 264989.88
$ws.Range("I141").Value = 149999
$ws.Range("J141").Value = 356982.6
$ws.Range("K141").Value = 149999
$ws.Range("L141").Value = 356982.6
$ws.Range("M141").Value = -144819
$ws.Range("N141").Value = -367342.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1616.3334
$ws.Range("I5").Value = 2300
$ws.Range("J5").Value = 932.6667
$ws.Range("K5").Value = 6900
$ws.Range("L5").Value = 2798.0001
$ws.Range("M5").Value = -6788
$ws.Range("N5").Value = -3022.0001

$ws.Range("H12").Value = 32.642857
$ws.Range("I12").Value = 45.4
$ws.Range("J12").Value = 25.555555
$ws.Range("K12").Value = 136.2
$ws.Range("L12").Value = 76.66666499999999
$ws.Range("M12").Value = 36.80000000000001
$ws.Range("N12").Value = -422.666665

$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("M24").ClearContents()

$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("M41").ClearContents()

$ws.Range("H109").Value = 3803.4
$ws.Range("I109").Value = 3803.4
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 11410.2
$ws.Range("L109").Value = 0
$ws.Range("M109").Value = -10370.2

$ws.Range("H116").Value = 858
$ws.Range("I116").Value = 858
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 2574
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 868

$ws.Range("H119").Value = 1500
$ws.Range("I119").Value = 1500
$ws.Range("J119").Value = 0
$ws.Range("K119").Value = 4500
$ws.Range("L119").Value = 0
$ws.Range("M119").Value = 338

$ws.Range("H132").Value = 1466.3334
$ws.Range("I132").Value = 400
$ws.Range("J132").Value = 1999.5
$ws.Range("K132").Value = 3600
$ws.Range("L132").Value = 17995.5
$ws.Range("M132").Value = -1070
$ws.Range("N132").Value = -23055.5

$ws.Range("H135").Value = 1616.3334
$ws.Range("I135").Value = 2300
$ws.Range("J135").Value = 932.6667
$ws.Range("K135").Value = 20700
$ws.Range("L135").Value = 8394.0003
$ws.Range("M135").Value = -18165
$ws.Range("N135").Value = -13464.0003

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 1000
$ws.Range("I13").Value = 1000
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 1000
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -860

$ws.Range("H61").Value = 2000
$ws.Range("I61").Value = 2000
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2000
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1798

$ws.Range("H113").Value = 2000
$ws.Range("I113").Value = 2000
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2000
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 170

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()

$ws.Range("H126").Value = 4820.8
$ws.Range("I126").Value = 5026
$ws.Range("J126").Value = 4000
$ws.Range("K126").Value = 15078
$ws.Range("L126").Value = 12000
$ws.Range("M126").Value = -12608
$ws.Range("N126").Value = -16940
